# Generate Report for Handoff
#
# The "b.md" file has now been handed off for translation. Update the
# Overview sheet's summary row for b.md, plus the per-locale detail
# sheets (zh-cn, de-de), to reflect the new handoff status/file/date.
#
# Cells are written in the same order the values first appear when the
# workbook is scanned sheet-by-sheet / row-by-row / column-by-column, so
# that newly introduced shared strings land in the same relative order
# as the canonical export.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet: row 3 is the "b.md" summary row ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-25 09:53:12"

# ---- zh-cn detail sheet: row 3 is the "b.md" detail row ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-25 09:53:03"

# ---- de-de detail sheet: row 3 is the "b.md" detail row ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-25 09:53:12"
